$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 88
$ws.Range("A88").Value = "What the maximum number of data files I can load?"
$ws.Range("B88").Value = "llama3.2:latest"
$ws.Range("C88").Value = "The maximum number of data files you can load is unlimited."

# Row 89
$ws.Range("A89").Value = "How many curves can I load in one go?"
$ws.Range("B89").Value = "llama3.2:latest"
$ws.Range("C89").Value = "You can load up to 450 curves at a time."

# Row 90
$ws.Range("A90").Value = "How many symbols can I have in the plot at any one time?"
$ws.Range("B90").Value = "llama3.2:latest"
$ws.Range("C90").Value = "You can have up to 10,000 symbols in a plot at any given time."
